$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column before I so everything from old column I ("견적서경로")
#    onward shifts one column to the right, making room for the new
#    "지급상태" (payment status) column while H keeps its position.
$ws.Columns("I:I").Insert()

# 2. Rename header H1 "Status" -> "입고상태" (delivery status) and fill the
#    newly inserted I1 with "지급상태" (payment status).
$ws.Range("H1").Value = "입고상태"
$ws.Range("I1").Value = "지급상태"

# 3. Remove the now-shifted trailing columns (old R..AG, now S..AH) that are
#    no longer part of the sheet.
$ws.Range("S1:AH1").EntireColumn.Delete()

# 4. Remove the second data row (old row 3) entirely.
$ws.Range("A3:AG3").EntireRow.Delete()

# 5. Update row 2 values to match the new layout / content.
#    F2 becomes a text "10" (was numeric) - pre-format as text so the
#    assigned value is stored as a string, then clear the format override
#    back to General so no stray number-format style lingers on the cell.
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "10"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "ㅁㄴㅇㄹㄴㅁㅇㄹㄷㅈㄹ"
$ws.Range("H2").Value = "미입고"
$ws.Range("I2").Value = "미지급"
# J2 (the shifted former I2) already holds an empty string from the column
# insert above - leave it untouched so the cell stays present-but-empty.
$ws.Range("K2").Value = "ㅁㄴㅇㄻㄴㅇ"
$ws.Range("L2").Value = "ㄻㄴㅇㄹㅁㄴㅇㄻㄴㅇㄹ"
$ws.Range("M2").Value = "ㅁㄴㅇㄻㄴㅇㄹ"
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 10000
$ws.Range("P2").Value = 10000
$ws.Range("Q2").Value = 1000
$ws.Range("R2").Value = 11000
